$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2021-10-08"

# Row 8 (June) - 2021 columns (T/U/V)
$ws.Range("T8").Value = 3
$ws.Range("U8").Value = 127
$ws.Range("V8").Value = 0.0231

# Row 12 (October) - update label and several year columns
$ws.Range("A12").Value = "October (through 10-08)"
$ws.Range("F12").Value = 15
$ws.Range("I12").Value = 10
$ws.Range("J12").Value = 0.1667
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 18
$ws.Range("M12").Value = 0.0526
$ws.Range("M12").NumberFormat = "0.0%"
$ws.Range("O12").Value = 8
$ws.Range("U12").Value = 56

# Row 13 (Total) - recomputed totals across years
$ws.Range("F13").Value = 398
$ws.Range("G13").Value = 0.1036
$ws.Range("I13").Value = 587
$ws.Range("J13").Value = 0.0814
$ws.Range("K13").Value = 62
$ws.Range("L13").Value = 505
$ws.Range("M13").Value = 0.1093
$ws.Range("O13").Value = 387
$ws.Range("P13").Value = 0.1
$ws.Range("T13").Value = 79
$ws.Range("U13").Value = 1226
$ws.Range("V13").Value = 0.0605
